# A new daily price record (Flame Seedless, 2023-12-20) is inserted as row 90
# of the "Uva" sheet, pushing the existing rows 90-169 down to 91-170.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 90 — this shifts every row
# from 90..169 down by one (new rows 91..170) and keeps all their values,
# styles, etc. intact.
$ws.Rows(90).Insert()

# Populate the newly-inserted row 90 with the new record. The "fixed"
# columns (market/region/product identifiers, quality, origin) repeat the
# same values used throughout this sheet.
$ws.Range("A90").Value = 8
$ws.Range("B90").Value = "Terminal La Palmera de La Serena"
$ws.Range("C90").Value = "Coquimbo"
$ws.Range("D90").Value = 45280
$ws.Range("E90").Value = 4
$ws.Range("F90").Value = "Fruta"
$ws.Range("G90").Value = 100109
$ws.Range("H90").Value = "Uva"
$ws.Range("I90").Value = 100109001
$ws.Range("J90").Value = "Uva"
$ws.Range("K90").Value = "Flame Seedless"
$ws.Range("L90").Value = "Primera"
$ws.Range("M90").Value = 520
$ws.Range("N90").Value = 9000
$ws.Range("O90").Value = 10000
$ws.Range("P90").Value = 9500
$ws.Range("Q90").Value = "`$/bandeja 10 kilos"
$ws.Range("R90").Value = "Provincia del Elqu$([char]0xED)"
$ws.Range("S90").Value = 950
$ws.Range("T90").Value = 10
